# Auto-generated Excel COM-interop script to apply the BuildingData fuel-type edits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BuildingData")

$ws.Range("E3").Value = "NGas"
$ws.Range("G3").Value = "NGas"
$ws.Range("E4").Value = "NGas"
$ws.Range("G4").Value = "NGas"
$ws.Range("E6").Value = "NGas"
$ws.Range("G6").Value = "NGas"
$ws.Range("E7").Value = "NGas"
$ws.Range("G7").Value = "NGas"
$ws.Range("E8").Value = "NGas"
$ws.Range("G8").Value = "NGas"
$ws.Range("E9").Value = "NGas"
$ws.Range("G9").Value = "NGas"
$ws.Range("E10").Value = "NGas"
$ws.Range("G10").Value = "NGas"
$ws.Range("E11").Value = "NGas"
$ws.Range("G11").Value = "NGas"
$ws.Range("E12").Value = "NGas"
$ws.Range("G12").Value = "NGas"
$ws.Range("E13").Value = "NGas"
$ws.Range("G13").Value = "NGas"
$ws.Range("E14").Value = "NGas"
$ws.Range("G14").Value = "NGas"
$ws.Range("E15").Value = "NGas"
$ws.Range("G15").Value = "NGas"
$ws.Range("E16").Value = "NGas"
$ws.Range("G16").Value = "NGas"
$ws.Range("E17").Value = "NGas"
$ws.Range("G17").Value = "NGas"
$ws.Range("E18").Value = "NGas"
$ws.Range("G18").Value = "NGas"
$ws.Range("E19").Value = "NGas"
$ws.Range("G19").Value = "NGas"
$ws.Range("E20").Value = "NGas"
$ws.Range("G20").Value = "NGas"
$ws.Range("E21").Value = "NGas"
$ws.Range("G21").Value = "NGas"
$ws.Range("E22").Value = "NGas"
$ws.Range("G22").Value = "NGas"
$ws.Range("E23").Value = "Oil4"
$ws.Range("G23").Value = "Oil4"
$ws.Range("E24").Value = "NGas"
$ws.Range("G24").Value = "NGas"
$ws.Range("E25").Value = "NGas"
$ws.Range("G25").Value = "NGas"
$ws.Range("E26").Value = "NGas"
$ws.Range("G26").Value = "NGas"
$ws.Range("E27").Value = "NGas"
$ws.Range("G27").Value = "NGas"
$ws.Range("E29").Value = "NGas"
$ws.Range("G29").Value = "NGas"
$ws.Range("E30").Value = "NGas"
$ws.Range("G30").Value = "NGas"
$ws.Range("E31").Value = "NGas"
$ws.Range("G31").Value = "NGas"
$ws.Range("E32").Value = "NGas"
$ws.Range("G32").Value = "NGas"
$ws.Range("E33").Value = "Steam"
$ws.Range("G33").Value = "Steam"
$ws.Range("E34").Value = "NGas"
$ws.Range("G34").Value = "NGas"
$ws.Range("E35").Value = "NGas"
$ws.Range("G35").Value = "NGas"
$ws.Range("E36").Value = "NGas"
$ws.Range("G36").Value = "NGas"
$ws.Range("E37").Value = "Steam"
$ws.Range("G37").Value = "Steam"
$ws.Range("E40").Value = "NGas"
$ws.Range("G40").Value = "NGas"
$ws.Range("E41").Value = "NGas"
$ws.Range("G41").Value = "NGas"
$ws.Range("E42").Value = "NGas"
$ws.Range("G42").Value = "NGas"
$ws.Range("E43").Value = "NGas"
$ws.Range("G43").Value = "NGas"
$ws.Range("E44").Value = "NGas"
$ws.Range("G44").Value = "NGas"
$ws.Range("E45").Value = "NGas"
$ws.Range("G45").Value = "NGas"
$ws.Range("E46").Value = "NGas"
$ws.Range("G46").Value = "NGas"
$ws.Range("E47").Value = "NGas"
$ws.Range("G47").Value = "NGas"
$ws.Range("E48").Value = "Steam"
$ws.Range("F48").Value = "Steam"
$ws.Range("G48").Value = "Steam"
$ws.Range("E49").Value = "Steam"
$ws.Range("G49").Value = "Steam"
$ws.Range("E50").Value = "Oil2"
$ws.Range("G50").Value = "Oil2"
$ws.Range("E51").Value = "Steam"
$ws.Range("G51").Value = "Steam"
$ws.Range("E52").Value = "NGas"
$ws.Range("G52").Value = "NGas"
$ws.Range("E53").Value = "NGas"
$ws.Range("G53").Value = "NGas"
$ws.Range("E54").Value = "NGas"
$ws.Range("G54").Value = "NGas"
$ws.Range("E55").Value = "NGas"
$ws.Range("G55").Value = "NGas"
$ws.Range("E56").Value = "Steam"
$ws.Range("F56").Value = "Steam"
$ws.Range("G56").Value = "Steam"
$ws.Range("E57").Value = "NGas"
$ws.Range("G57").Value = "NGas"
$ws.Range("E58").Value = "NGas"
$ws.Range("G58").Value = "NGas"
$ws.Range("E59").Value = "NGas"
$ws.Range("G59").Value = "NGas"
$ws.Range("E60").Value = "Steam"
$ws.Range("F60").Value = "Steam"
$ws.Range("G60").Value = "Steam"
$ws.Range("E61").Value = "NGas"
$ws.Range("G61").Value = "NGas"
$ws.Range("E62").Value = "NGas"
$ws.Range("G62").Value = "NGas"
$ws.Range("E63").Value = "NGas"
$ws.Range("G63").Value = "NGas"
$ws.Range("E65").Value = "NGas"
$ws.Range("G65").Value = "NGas"
$ws.Range("E66").Value = "NGas"
$ws.Range("G66").Value = "NGas"
$ws.Range("E67").Value = "NGas"
$ws.Range("G67").Value = "NGas"
$ws.Range("E68").Value = "NGas"
$ws.Range("G68").Value = "NGas"

# Update the view/selection to match the saved state (scrolled down, new active cell)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 51
$ws.Range("E69").Select()
